$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 3 new daily rows (06, 07, 08-Oct-2021) to the currency-parity table.
# Column A holds a date-like label but is stored as plain text in the source
# data, so we force a Text format for the assignment (otherwise the COM layer
# auto-converts "DD-MM-YYYY"-looking strings into real date serials) and then
# clear the format back off so the cell is indistinguishable from its neighbours.

# Row 193 (06-10-2021)
$ws.Cells.Item(193, 1).NumberFormat = "@"
$ws.Cells.Item(193, 1).Value = "06-10-2021"
$ws.Cells.Item(193, 1).ClearFormats()
$ws.Cells.Item(193, 2).Value = 33.812
$ws.Cells.Item(193, 3).Value = 1
$ws.Cells.Item(193, 4).Value = 4.1651
$ws.Cells.Item(193, 5).Value = 6.86
$ws.Cells.Item(193, 6).Value = 626.55
$ws.Cells.Item(193, 7).Value = 21.838
$ws.Cells.Item(193, 8).Value = 6.41
$ws.Cells.Item(193, 9).Value = 128.13
$ws.Cells.Item(193, 10).Value = 8.527900000000001
$ws.Cells.Item(193, 11).Value = 8.7301
$ws.Cells.Item(193, 12).Value = 0.7079
$ws.Cells.Item(193, 13).Value = 3.672
$ws.Cells.Item(193, 14).Value = 9.0403
$ws.Cells.Item(193, 15).Value = 1.3706
$ws.Cells.Item(193, 16).Value = 1.2571
$ws.Cells.Item(193, 17).Value = 1
$ws.Cells.Item(193, 18).Value = 2.103
$ws.Cells.Item(193, 19).Value = 0.82
$ws.Cells.Item(193, 20).Value = 1
$ws.Cells.Item(193, 21).Value = 1.3566
$ws.Cells.Item(193, 22).Value = 7.7848
$ws.Cells.Item(193, 23).Value = 1.4349
$ws.Cells.Item(193, 24).Value = 27.885
$ws.Cells.Item(193, 25).Value = 22756
$ws.Cells.Item(193, 26).Value = 0.8618
$ws.Cells.Item(193, 27).Value = 308.39
$ws.Cells.Item(193, 28).Value = 102.3612
$ws.Cells.Item(193, 29).Value = 0.9277
$ws.Cells.Item(193, 30).Value = 6920.5
$ws.Cells.Item(193, 31).Value = 26.3965
$ws.Cells.Item(193, 32).Value = 4.26
$ws.Cells.Item(193, 33).Value = 15.6697
$ws.Cells.Item(193, 34).Value = 0.7334000000000001
$ws.Cells.Item(193, 35).Value = 8.871
$ws.Cells.Item(193, 36).Value = 4.1285
$ws.Cells.Item(193, 37).Value = 98.8875
$ws.Cells.Item(193, 38).Value = 810.63
$ws.Cells.Item(193, 39).Value = 3784.05
$ws.Cells.Item(193, 40).Value = 24
$ws.Cells.Item(193, 41).Value = 56.161
$ws.Cells.Item(193, 42).Value = 50.643
$ws.Cells.Item(193, 43).Value = 20.5331
$ws.Cells.Item(193, 44).Value = 42.85
$ws.Cells.Item(193, 45).Value = 7.7272
$ws.Cells.Item(193, 46).Value = 14.9867
$ws.Cells.Item(193, 47).Value = 5.4681
$ws.Cells.Item(193, 48).Value = 42000
$ws.Cells.Item(193, 49).Value = 3.75
$ws.Cells.Item(193, 50).Value = 4.1785
$ws.Cells.Item(193, 51).Value = 72.28830000000001
$ws.Cells.Item(193, 52).Value = 74.44750000000001
$ws.Cells.Item(193, 53).Value = 14250
$ws.Cells.Item(193, 54).Value = 170.5
$ws.Cells.Item(193, 55).Value = 3.2306
$ws.Cells.Item(193, 56).Value = 424.7
$ws.Cells.Item(193, 57).Value = 1188.7
$ws.Cells.Item(193, 58).Value = 111.45
$ws.Cells.Item(193, 59).Value = 6.4454
$ws.Cells.Item(193, 60).Value = 3.9732

# Row 194 (07-10-2021)
$ws.Cells.Item(194, 1).NumberFormat = "@"
$ws.Cells.Item(194, 1).Value = "07-10-2021"
$ws.Cells.Item(194, 1).ClearFormats()
$ws.Cells.Item(194, 2).Value = 33.817
$ws.Cells.Item(194, 3).Value = 1
$ws.Cells.Item(194, 4).Value = 4.184
$ws.Cells.Item(194, 5).Value = 6.86
$ws.Cells.Item(194, 6).Value = 626.8099999999999
$ws.Cells.Item(194, 7).Value = 22.001
$ws.Cells.Item(194, 8).Value = 6.4419
$ws.Cells.Item(194, 9).Value = 128.58
$ws.Cells.Item(194, 10).Value = 8.5916
$ws.Cells.Item(194, 11).Value = 8.8094
$ws.Cells.Item(194, 12).Value = 0.7094
$ws.Cells.Item(194, 13).Value = 3.672
$ws.Cells.Item(194, 14).Value = 9.074999999999999
$ws.Cells.Item(194, 15).Value = 1.3776
$ws.Cells.Item(194, 16).Value = 1.2604
$ws.Cells.Item(194, 17).Value = 1
$ws.Cells.Item(194, 18).Value = 2.103
$ws.Cells.Item(194, 19).Value = 0.82
$ws.Cells.Item(194, 20).Value = 1
$ws.Cells.Item(194, 21).Value = 1.3599
$ws.Cells.Item(194, 22).Value = 7.7866
$ws.Cells.Item(194, 23).Value = 1.4476
$ws.Cells.Item(194, 24).Value = 27.977
$ws.Cells.Item(194, 25).Value = 22759
$ws.Cells.Item(194, 26).Value = 0.866
$ws.Cells.Item(194, 27).Value = 310.79
$ws.Cells.Item(194, 28).Value = 103.3354
$ws.Cells.Item(194, 29).Value = 0.9282
$ws.Cells.Item(194, 30).Value = 6912.6
$ws.Cells.Item(194, 31).Value = 26.335
$ws.Cells.Item(194, 32).Value = 4.2809
$ws.Cells.Item(194, 33).Value = 15.6538
$ws.Cells.Item(194, 34).Value = 0.7372
$ws.Cells.Item(194, 35).Value = 8.8773
$ws.Cells.Item(194, 36).Value = 4.1345
$ws.Cells.Item(194, 37).Value = 98.9126
$ws.Cells.Item(194, 38).Value = 816.28
$ws.Cells.Item(194, 39).Value = 3786.01
$ws.Cells.Item(194, 40).Value = 24
$ws.Cells.Item(194, 41).Value = 56.152
$ws.Cells.Item(194, 42).Value = 50.865
$ws.Cells.Item(194, 43).Value = 20.6459
$ws.Cells.Item(194, 44).Value = 42.8901
$ws.Cells.Item(194, 45).Value = 7.7272
$ws.Cells.Item(194, 46).Value = 15.0508
$ws.Cells.Item(194, 47).Value = 5.5214
$ws.Cells.Item(194, 48).Value = 42000
$ws.Cells.Item(194, 49).Value = 3.75
$ws.Cells.Item(194, 50).Value = 4.1825
$ws.Cells.Item(194, 51).Value = 72.43689999999999
$ws.Cells.Item(194, 52).Value = 74.98
$ws.Cells.Item(194, 53).Value = 14250
$ws.Cells.Item(194, 54).Value = 170.75
$ws.Cells.Item(194, 55).Value = 3.2369
$ws.Cells.Item(194, 56).Value = 424.68
$ws.Cells.Item(194, 57).Value = 1192.45
$ws.Cells.Item(194, 58).Value = 111.43
$ws.Cells.Item(194, 59).Value = 6.4576
$ws.Cells.Item(194, 60).Value = 3.946

# Row 195 (08-10-2021)
$ws.Cells.Item(195, 1).NumberFormat = "@"
$ws.Cells.Item(195, 1).Value = "08-10-2021"
$ws.Cells.Item(195, 1).ClearFormats()
$ws.Cells.Item(195, 2).Value = 33.773
$ws.Cells.Item(195, 3).Value = 1
$ws.Cells.Item(195, 4).Value = 4.1551
$ws.Cells.Item(195, 5).Value = 6.86
$ws.Cells.Item(195, 6).Value = 627.2
$ws.Cells.Item(195, 7).Value = 22.0013
$ws.Cells.Item(195, 8).Value = 6.4348
$ws.Cells.Item(195, 9).Value = 128.77
$ws.Cells.Item(195, 10).Value = 8.5616
$ws.Cells.Item(195, 11).Value = 8.7759
$ws.Cells.Item(195, 12).Value = 0.7088
$ws.Cells.Item(195, 13).Value = 3.672
$ws.Cells.Item(195, 14).Value = 9.048
$ws.Cells.Item(195, 15).Value = 1.3667
$ws.Cells.Item(195, 16).Value = 1.2553
$ws.Cells.Item(195, 17).Value = 1
$ws.Cells.Item(195, 18).Value = 2.1079
$ws.Cells.Item(195, 19).Value = 0.82
$ws.Cells.Item(195, 20).Value = 1
$ws.Cells.Item(195, 21).Value = 1.3581
$ws.Cells.Item(195, 22).Value = 7.7842
$ws.Cells.Item(195, 23).Value = 1.4413
$ws.Cells.Item(195, 24).Value = 27.96
$ws.Cells.Item(195, 25).Value = 22758
$ws.Cells.Item(195, 26).Value = 0.8648
$ws.Cells.Item(195, 27).Value = 309.8
$ws.Cells.Item(195, 28).Value = 102.6367
$ws.Cells.Item(195, 29).Value = 0.9277
$ws.Cells.Item(195, 30).Value = 6901
$ws.Cells.Item(195, 31).Value = 26.3
$ws.Cells.Item(195, 32).Value = 4.2772
$ws.Cells.Item(195, 33).Value = 15.66
$ws.Cells.Item(195, 34).Value = 0.7341
$ws.Cells.Item(195, 35).Value = 8.871600000000001
$ws.Cells.Item(195, 36).Value = 4.0915
$ws.Cells.Item(195, 37).Value = 98.94
$ws.Cells.Item(195, 38).Value = 813.62
$ws.Cells.Item(195, 39).Value = 3777.76
$ws.Cells.Item(195, 40).Value = 24
$ws.Cells.Item(195, 41).Value = 56.221
$ws.Cells.Item(195, 42).Value = 50.56
$ws.Cells.Item(195, 43).Value = 20.619
$ws.Cells.Item(195, 44).Value = 43.19
$ws.Cells.Item(195, 45).Value = 7.7054
$ws.Cells.Item(195, 46).Value = 14.9242
$ws.Cells.Item(195, 47).Value = 5.5014
$ws.Cells.Item(195, 48).Value = 42000
$ws.Cells.Item(195, 49).Value = 3.75
$ws.Cells.Item(195, 50).Value = 4.182
$ws.Cells.Item(195, 51).Value = 71.7338
$ws.Cells.Item(195, 52).Value = 74.78
$ws.Cells.Item(195, 53).Value = 14215
$ws.Cells.Item(195, 54).Value = 170.25
$ws.Cells.Item(195, 55).Value = 3.2238
$ws.Cells.Item(195, 56).Value = 425.45
$ws.Cells.Item(195, 57).Value = 1190.57
$ws.Cells.Item(195, 58).Value = 111.58
$ws.Cells.Item(195, 59).Value = 6.4533
$ws.Cells.Item(195, 60).Value = 3.9566
